$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.135.35"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.573.41"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'207.23"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'22.27"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.798.34"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.578.46"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "27.143.04"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "'7.37"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'214.17"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "0.0₃0683"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  -3.41%  "
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "'152.59"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("D27").Value = "'14.92"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").Value = "'0.103"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").Value = "1.395.40"
$ws.Range("E33").Value = "  +2.43%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "'0.946"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'0.812"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +3.32%  "
$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("D44").Value = "'5.41"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "'2.20"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").Value = "'63.76"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "1.710.74"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("D48").Value = "'85.68"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "0.0₇0990"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "'0.0953"
$ws.Range("E51").Value = "  -0.39%  "
